$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Company")
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(3).Delete()

$ws2 = $wb.Worksheets.Item("AddCoverageTeam")
$ws2.Range("B2").Value = "Jacklyn Robinson"
$ws2.Range("C11").Select()
$ws2.Activate()

Write-Host "done"
